# alteracao nas classes GenericLogic, GenericStep, LoginLogic e LoginPage
# GenericLogic gera um novo usuario de teste aleatorio (ex: "teste" + letras)
# e LoginLogic/LoginPage registram o valor utilizado na planilha de controle.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuario")

$ws.Range("A2").Value = "testejtytnm"
